# Restored from revision of admin on 10/29/2020 07:59:55 AM.TEST
# Update rule R30's "From" threshold (cell C10 on the "Rules" sheet) from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
